# ezanaLMS_Students.xlsx: add a new "M" column populated with a constant
# hash-like token on every data row (header included) - this is the extra
# bind-variable column that mysqli_stmt::bind_param() now expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$token = "53c904468e7edec9a7f2501d8a8c8d5140c434cb"

# Header + all 14 data rows (sheet has rows 1..15) get the same value.
for ($r = 1; $r -le 15; $r++) {
    $cell = $ws.Cells.Item($r, 13)
    $cell.Value = $token
    # Give the new column its own (non-default) plain-black style, distinct
    # from both the workbook default and the blue hyperlink style used by D3.
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.Color = 0
}

# Move the active selection to the last cell touched, like a user who just
# finished typing the new column would leave it.
$ws.Range("M15").Select()
